# Rename the embedded logo pictures' shape names.
#
# The document has two headers (first-page + default, both carrying the
# BTec_Logo-Orange JPEG) and two footers (first-page + default, both
# carrying the Pearson "PearsonLogo.png" PNG). Word stores a default
# "shape name" for every picture (visible in the Selection Pane / the
# picture's docPr@name) separate from its alt-text description. This
# edit renumbers those default picture names:
#   - the BTec jpg logo: image1.jpg -> image2.jpg (in both headers)
#   - the Pearson png logo: image2.png -> image1.png (in both footers)

$d = $word.ActiveDocument

function Rename-HeaderFooterPicture($story, $newName) {
    if ($story.Exists -and $story.Range.InlineShapes.Count -gt 0) {
        for ($i = 1; $i -le $story.Range.InlineShapes.Count; $i++) {
            $shp = $story.Range.InlineShapes.Item($i)
            $shp.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($hfIndex = 1; $hfIndex -le 3; $hfIndex++) {
        $hdr = $sec.Headers.Item($hfIndex)
        Rename-HeaderFooterPicture $hdr "image2.jpg"

        $ftr = $sec.Footers.Item($hfIndex)
        Rename-HeaderFooterPicture $ftr "image1.png"
    }
}

Write-Output "done"
